$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.67"
$ws.Range("E2").Value = "'6.37%"
$ws.Range("D3").Value = "'40.54"
$ws.Range("E3").Value = "'11.75%"
$ws.Range("D4").Value = "'5.928"
$ws.Range("E4").Value = "'15.42%"
$ws.Range("D5").Value = "'0.08136"
$ws.Range("E5").Value = "'5.30%"
$ws.Range("D6").Value = "'4.578"
$ws.Range("E6").Value = "'4.27%"
$ws.Range("E7").Value = "'5.33%"
$ws.Range("D8").Value = "'1.954"
$ws.Range("E8").Value = "'5.36%"
$ws.Range("E9").Value = "'0.00%"
$ws.Range("D10").Value = "'0.9445"
$ws.Range("E10").Value = "'2.38%"
$ws.Range("D11").Value = "'0.1311"
$ws.Range("E11").Value = "'16.64%"
$ws.Range("D12").Value = "'0.1995"
$ws.Range("E12").Value = "'7.33%"
$ws.Range("D13").Value = "'0.09266"
$ws.Range("E13").Value = "'5.26%"
$ws.Range("D14").Value = "'0.03420"
$ws.Range("E14").Value = "'2.90%"
$ws.Range("E15").Value = "'1.04%"
$ws.Range("D16").Value = "'0.001340"
$ws.Range("E16").Value = "'-3.26%"
$ws.Range("D17").Value = "'0.005995"
$ws.Range("E17").Value = "'-1.66%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.373"
$ws.Range("E18").Value = "'0.14%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3496"
$ws.Range("E19").Value = "'1.43%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'7.738"
$ws.Range("E20").Value = "'22.31%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1440"
$ws.Range("E21").Value = "'10.84%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2448"
$ws.Range("E22").Value = "'5.79%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04433"
$ws.Range("E23").Value = "'2.15%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001254"
$ws.Range("E24").Value = "'4.23%"
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").Value = "'0.004355"
$ws.Range("E25").Value = "'2.15%"
$ws.Range("D26").Value = "'0.0001190"
$ws.Range("E26").Value = "'-10.55%"
$ws.Range("D27").Value = "'0.0003992"
$ws.Range("E27").Value = "'37.46%"
$ws.Range("D39").Value = "'0.02498"
$ws.Range("E39").Value = "'19.05%"
$ws.Range("D40").Value = "'0.05282"
$ws.Range("E40").Value = "'7.19%"
$ws.Range("D41").Value = "'0.007604"
$ws.Range("D42").Value = "'0.1432"
$ws.Range("E42").Value = "'6.20%"
$ws.Range("D43").Value = "'0.008949"
$ws.Range("E43").Value = "'5.40%"
$ws.Range("D44").Value = "'0.002065"
$ws.Range("E44").Value = "'-0.39%"
$ws.Range("D45").Value = "'0.009491"
$ws.Range("E45").Value = "'12.99%"
$ws.Range("D46").Value = "'0.00006882"
$ws.Range("E46").Value = "'6.25%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.05%"
$ws.Range("D48").Value = "'0.002897"
$ws.Range("E48").Value = "'-12.09%"
$ws.Range("D49").Value = "'0.001801"
$ws.Range("E49").Value = "'24.64%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.05%"
